$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.100.63"
$ws.Range("E2").Value = "  -1.71%  "

$ws.Range("D3").Value = "3.512.53"
$ws.Range("E3").Value = "  -0.61%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'576.43"
$ws.Range("E5").Value = "  +2.96%  "

$ws.Range("D6").Value = "'179.10"
$ws.Range("E6").Value = "  -7.35%  "

$ws.Range("D7").Value = "'0.634"
$ws.Range("E7").Value = "  +3.55%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").Value = "'0.633"
$ws.Range("E9").Value = "  -1.69%  "

$ws.Range("D10").Value = "'0.157"
$ws.Range("E10").Value = "  +3.32%  "

$ws.Range("D11").Value = "'54.69"
$ws.Range("E11").Value = "  -3.88%  "

$ws.Range("D12").Value = "'0.0000273"
$ws.Range("E12").Value = "  -0.03%  "

$ws.Range("D13").Value = "'9.22"
$ws.Range("E13").Value = "  -3.39%  "

$ws.Range("D14").Value = "4.060.05"
$ws.Range("E14").Value = "  -1.04%  "

$ws.Range("D15").Value = "3.504.86"
$ws.Range("E15").Value = "  -0.87%  "

$ws.Range("E16").Value = "  -0.28%  "

$ws.Range("D17").Value = "'18.33"
$ws.Range("E17").Value = "  -0.81%  "

$ws.Range("D18").Value = "'12.13"
$ws.Range("E18").Value = "  +1.12%  "

$ws.Range("D19").Value = "66.054.15"
$ws.Range("E19").Value = "  -1.77%  "

$ws.Range("D20").Value = "'1.01"
$ws.Range("E20").Value = "  -0.23%  "

$ws.Range("D21").Value = "'415.30"
$ws.Range("E21").Value = "  +1.53%  "

$ws.Range("D22").Value = "'4.20"
$ws.Range("E22").Value = "  +5.02%  "

$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").Value = "'85.37"
$ws.Range("E23").Value = "  -0.79%  "

$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").Value = "'4.25"
$ws.Range("E24").Value = "  -0.25%  "

$ws.Range("D25").Value = "'12.74"
$ws.Range("E25").Value = "  +5.26%  "

$ws.Range("D26").Value = "'10.91"
$ws.Range("E26").Value = "  -3.25%  "

$ws.Range("D27").Value = "'2.85"
$ws.Range("E27").Value = "  -4.03%  "

$ws.Range("D28").Value = "'9.00"
$ws.Range("E28").Value = "  +0.49%  "

$ws.Range("D29").Value = "'30.32"
$ws.Range("E29").Value = "  -1.25%  "

$ws.Range("D30").Value = "'619.65"
$ws.Range("E30").Value = "  -9.18%  "

$ws.Range("D31").Value = "'6.42"
$ws.Range("E31").Value = "  -5.85%  "

$ws.Range("D32").Value = "'11.65"
$ws.Range("E32").Value = "  -1.94%  "

$ws.Range("D33").Value = "'0.110"
$ws.Range("E33").Value = "  -2.37%  "

$ws.Range("D34").Value = "'59.71"
$ws.Range("E34").Value = "  -1.65%  "

$ws.Range("D35").Value = "'0.152"
$ws.Range("E35").Value = "  +9.36%  "

$ws.Range("B36").Value = "PEPE"
$ws.Range("C36").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D36").Value = "0.0₃0806"
$ws.Range("E36").Value = "  -3.00%  "

$ws.Range("B37").Value = "Dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  +0.10%  "

$ws.Range("D38").Value = "'37.41"
$ws.Range("E38").Value = "  -4.95%  "

$ws.Range("D39").Value = "3.287.63"
$ws.Range("E39").Value = "  +8.43%  "

$ws.Range("D40").Value = "'0.380"
$ws.Range("E40").Value = "  -4.56%  "

$ws.Range("D41").Value = "'3.34"
$ws.Range("E41").Value = "  -1.81%  "

$ws.Range("D42").Value = "'0.997"
$ws.Range("E42").Value = "  -0.20%  "

$ws.Range("D43").Value = "'2.92"
$ws.Range("E43").Value = "  -4.38%  "

$ws.Range("D44").Value = "'0.0418"
$ws.Range("E44").Value = "  -0.80%  "

$ws.Range("D45").Value = "'2.51"
$ws.Range("E45").Value = "  -6.56%  "

$ws.Range("D46").Value = "'3.24"
$ws.Range("E46").Value = "  -3.51%  "

$ws.Range("D47").Value = "'2.71"
$ws.Range("E47").Value = "  -0.96%  "

$ws.Range("D48").Value = "'0.133"
$ws.Range("E48").Value = "  +0.91%  "

$ws.Range("D49").Value = "'138.60"
$ws.Range("E49").Value = "  -0.62%  "

$ws.Range("D50").Value = "'8.46"
$ws.Range("E50").Value = "  -7.83%  "

$ws.Range("D51").Value = "'2.30"
$ws.Range("E51").Value = "  -10.16%  "

Write-Output "done"